$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "VFQA_Test652"
$ws.Range("D4").Value = "10121360348"
$ws.Range("D5").Value = "1-4NDZVP1"
$ws.Range("D6").Value = "1-10121360466"
$ws.Range("D7").Value = "24-1-2018"
$ws.Range("D9").Value = "97478152659"
$ws.Range("D10").Value = "8962702800901163659"
